$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new data rows just above the old blank separator / summary block
# (rows 67-68), pushing the separator + summary rows down by two.
$ws.Rows("67:68").Insert()

# New data row 67 (2014-03-15, 12:10 -> 13:10, i.e. 60 min)
$ws.Range("A67").Value = 2014
$ws.Range("B67").Value = 3
$ws.Range("C67").Value = 15
$ws.Range("D67").Value = 0.50694444444444442
$ws.Range("E67").Value = 0.54861111111111105
$ws.Range("F67").Formula = "=(E67-D67)*24*60"
$ws.Range("G67").Formula = "=F67/60"

# New data row 68 (2014-03-15, 14:15 -> 18:00, i.e. 225 min)
$ws.Range("A68").Value = 2014
$ws.Range("B68").Value = 3
$ws.Range("C68").Value = 15
$ws.Range("D68").Value = 0.59375
$ws.Range("E68").Value = 0.75
$ws.Range("F68").Formula = "=(E68-D68)*24*60"
$ws.Range("G68").Formula = "=F68/60"

# Move the view selection the way the author's workbook ended up (A69)
[void]$ws.Range("A69").Select()
